# EA 23.209 Stable: merchants recruited in tents are no longer temporary.
# Update the "title05" row (B24:D24) translation strings to drop "Temporary"/"臨時"/"临时".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B24").Value = "Merchant Recruitment"
$ws.Range("C24").Value = "商人の雇用"
$ws.Range("D24").Value = "商人招募"

# The teleporter toggle cells (C28:D29) were styled with a "Noto Sans SC" font
# variant that lacked a charset; re-apply the font so they share the same
# (charset-carrying) font entry used elsewhere in the sheet.
$ws.Range("C28:D29").Font.Name = "Noto Sans SC"
$ws.Range("C28:D29").Font.Size = 10

# Restore the last saved selection to D29 (bottom-right translated cell).
$ws.Range("D29").Select()
